# Figure and Table Edits
# Figure and Table edits for publication.
# Updates statistical results (estimates, p-values, etc.) in Sheet1 of the
# Stats Tables workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rows 16-21 (emmeans / contrast table) ---
$ws.Range("B16").Value = 0.0163
$ws.Range("D16").Value = 0.9
$ws.Range("E16").Value = 9.93
$ws.Range("G16").Value = "0.0016*"

$ws.Range("B17").Value = 0.396
$ws.Range("E17").Value = 31.4

$ws.Range("B18").Value = 0.0024
$ws.Range("D18").Value = 0.96

$ws.Range("B19").Value = 0.0063
$ws.Range("D19").Value = 0.94
$ws.Range("E19").Value = 3.67
$ws.Range("G19").Value = 0.055

$ws.Range("B20").Value = 0.0043
$ws.Range("D20").Value = 0.95
$ws.Range("E20").Value = 0.732
$ws.Range("G20").Value = 0.39

$ws.Range("B21").Value = 0.0023
$ws.Range("E21").Value = 4.61
$ws.Range("G21").Value = "0.032*"

# --- Rows 43-48 (contrast table, Microsite:Year) ---
$ws.Range("M43").Value = "1.88 (0.353)"
$ws.Range("N43").Value = 5.34

$ws.Range("M44").Value = "-2.84 (1.10)"
$ws.Range("N44").Value = -2.58
$ws.Range("O44").Value = "0.012*"

$ws.Range("M45").Value = "1.14 (0.521)"
$ws.Range("N45").Value = 2.1800000000000002
$ws.Range("O45").Value = "0.032*"

$ws.Range("M46").Value = "2.63 (0.467)"
$ws.Range("N46").Value = 5.63

$ws.Range("M47").Value = "-3.58 (1.08)"
$ws.Range("N47").Value = -3.33
$ws.Range("O47").Value = "0.0013*"

$ws.Range("M48").Value = "-2.09 (1.23)"
$ws.Range("N48").Value = -1.7
$ws.Range("O48").Value = 0.092

# --- Sheet view state: scroll position and active selection ---
$ws.Activate()
$ws.Range("A21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("P46").Select() | Out-Null
